$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for the Action-table entry "close workflow" above the
# existing "batch execute workflow" (A43) action row. This shifts the
# blank separator rows / resource header / resource table all down by one.
$ws.Rows.Item(43).Insert()

# Insert a new row for the matching Resource-table entry above the existing
# "batch execute workflow" (A61, now shifted to A62) resource row.
$ws.Rows.Item(62).Insert()

# Re-apply the formatting of sibling data rows to the two freshly inserted
# (blank/unformatted) rows so they look like the other table rows instead of
# the blank spacer row they inherited format from.
$ws.Range("A44:I44").Copy()
$ws.Range("A43:I43").PasteSpecial(-4122)

$ws.Range("A63:I63").Copy()
$ws.Range("A62:I62").PasteSpecial(-4122)

# --- Resource-table row (A62:G62): "close workflow" resource ---
$ws.Range("B62").Value = "关闭工作流"
$ws.Range("A62").Value = "0c235fe1-24fc-4db1-8be6-a131ff7dfd6c"
$ws.Range("C62").Value = "UNION"
$ws.Range("D62").Value = "res.up.flow.close"
$ws.Range("E62").Value = "w.todo"
$ws.Range("F62").Value = "resource.flow"
$ws.Range("G62").Value = 8

# --- Action-table row (A43:H43): "close workflow" action ---
$ws.Range("A43").Value = "849f2149-4abd-4213-9d47-0d1d0b12d051"
$ws.Range("B43").Formula = "=A62"
$ws.Range("C43").Formula = "=A25"
$ws.Range("D43").Value = "act.up.flow.close"
$ws.Range("E43").Value = "PUT"
$ws.Range("F43").Value = "/api/up/flow/close"
$ws.Range("G43").Value = "关闭工作流"
$ws.Range("H43").Value = 8

# Move the visible selection to D44 (matches the author's final cursor spot).
$ws.Range("D44").Select()

Write-Host "Workflow close action/resource rows inserted"
